# daily auto push: 2026-02-20 19:03 UTC
# Two new timestamped log entries (2026/02/20 23:00 and 2026/02/21 03:00)
# were appended to the "sei2" time-of-day log sheet. This pushes all the
# later rows down by two and appends two rows at the bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift every row from 839 downward by two rows, opening up space for the
# two freshly-logged entries (keeps all the later date/weekday/hour rows
# intact, just relocated).
$ws.Rows.Item(839).Resize(2).Insert()

# New rows to insert at 839-840.
$newRows = @(
    @("2026/02/20", "金", 23, 201),
    @("2026/02/21", "土", 3, 201)
)

for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = 839 + $i
    $entry = $newRows[$i]

    # Column A holds a literal "yyyy/mm/dd" string (not a real date value),
    # matching the rest of the sheet -- force Text format first so Excel
    # doesn't auto-convert it to a date serial, then restore the default
    # "Normal" style so no stray formatting is left behind.
    $dateCell = $ws.Cells.Item($r, 1)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $entry[0]
    $dateCell.Style = "Normal"

    $ws.Cells.Item($r, 2).Value = $entry[1]
    $ws.Cells.Item($r, 3).Value = $entry[2]
    $ws.Cells.Item($r, 4).Value = $entry[3]
}
